# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 22:34"

# Row 4 - Estados Unidos: refreshed case numbers
$ws.Range("B4").Value = 1341281
$ws.Range("C4").Value = 19496
$ws.Range("E4").Value = 1029098
$ws.Range("F4").Value = 16796
$ws.Range("G4").Value = 1208
$ws.Range("H4").Value = 79823

# Row 10 - Alemania: refreshed case numbers
$ws.Range("B10").Value = 171264
$ws.Range("C10").Value = 676
$ws.Range("E10").Value = 20421
$ws.Range("G10").Value = 33
$ws.Range("H10").Value = 7543

# Row 15 - Canada: refreshed case numbers
$ws.Range("B15").Value = 67643
$ws.Range("C15").Value = 1209
$ws.Range("D15").Value = 31065
$ws.Range("E15").Value = 31886
$ws.Range("G15").Value = 123
$ws.Range("H15").Value = 4692

# Rows 84-86: Costa de Marfil moves ahead of Senegal / Republica de Macedonia,
# each row keeping the numbers that follow the new country ordering.
$ws.Range("A84").Value = "Costa de Marfil"
$ws.Range("B84").Value = 1667
$ws.Range("C84").Value = 65
$ws.Range("D84").Value = 769
$ws.Range("E84").Value = 877
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 21

$ws.Range("A85").Value = "Senegal"
$ws.Range("B85").Value = 1634
$ws.Range("C85").Value = 83
$ws.Range("D85").Value = 643
$ws.Range("E85").Value = 974
$ws.Range("F85").Value = 6
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = 17

$ws.Range("A86").Value = "Republica de Macedonia"
$ws.Range("B86").Value = 1622
$ws.Range("C86").Value = 36
$ws.Range("D86").Value = 1112
$ws.Range("E86").Value = 419
$ws.Range("F86").Value = 21
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 91

# Rows 104-105: Niger moves ahead of Libano
$ws.Range("A104").Value = "Niger"
$ws.Range("B104").Value = 815
$ws.Range("C104").Value = 20
$ws.Range("D104").Value = 617
$ws.Range("E104").Value = 153
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 45

$ws.Range("A105").Value = "Libano"
$ws.Range("B105").Value = 809
$ws.Range("C105").Value = 13
$ws.Range("D105").Value = 234
$ws.Range("E105").Value = 549
$ws.Range("F105").Value = 4
$ws.Range("H105").Value = 26

# Rows 137-138: Ruanda moves ahead of Congo
$ws.Range("A137").Value = "Ruanda"
$ws.Range("B137").Value = 280
$ws.Range("C137").Value = 7
$ws.Range("D137").Value = 140
$ws.Range("E137").Value = 140
$ws.Range("H137").Value = 0

$ws.Range("A138").Value = "Congo"
$ws.Range("B138").Value = 274
$ws.Range("D138").Value = 33
$ws.Range("E138").Value = 231
$ws.Range("H138").Value = 10
